$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1) - rows 5-37
$wsA = $wb.Worksheets.Item("展览")
$wsA.Range("F5").Value = 1156
$wsA.Range("F6").Value = 14311
$wsA.Range("F7").Value = 16403
$wsA.Range("F9").Value = 94
$wsA.Range("F10").Value = 6
$wsA.Range("F21").Value = 1251
$wsA.Range("F24").Value = 36
$wsA.Range("F25").Value = 1
$wsA.Range("F26").Value = 6635
$wsA.Range("F27").Value = 970
$wsA.Range("F28").Value = 2
$wsA.Range("F29").Value = 18
$wsA.Range("F32").Value = 5727
$wsA.Range("F35").Value = 183
$wsA.Range("F36").Value = 4779
$wsA.Range("F37").Value = 17

# Sheet "全部类型" (sheetId 4) - rows 5-40 (one extra row vs 展览)
$wsB = $wb.Worksheets.Item("全部类型")
$wsB.Range("F5").Value = 1156
$wsB.Range("F6").Value = 14311
$wsB.Range("F7").Value = 16403
$wsB.Range("F9").Value = 94
$wsB.Range("F10").Value = 6
$wsB.Range("F21").Value = 1251
$wsB.Range("F25").Value = 36
$wsB.Range("F26").Value = 1
$wsB.Range("F27").Value = 6635
$wsB.Range("F28").Value = 970
$wsB.Range("F29").Value = 2
$wsB.Range("F30").Value = 18
$wsB.Range("F35").Value = 5727
$wsB.Range("F38").Value = 183
$wsB.Range("F39").Value = 4779
$wsB.Range("F40").Value = 17
